$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("H50").Value = 79.3
$ws.Range("G51").Value = 71.09999999999999
$ws.Range("D52").Value = 90.90000000000001
$ws.Range("D53").Value = 90.5
$ws.Range("E53").Value = 85.3
$ws.Range("B54").Value = 80.8
$ws.Range("D54").Value = 89.40000000000001
$ws.Range("D55").Value = 92.09999999999999
$ws.Range("E55").Value = 83.59999999999999
$ws.Range("H55").Value = 80.90000000000001
$ws.Range("C56").Value = 87.5
$ws.Range("D56").Value = 94.2
$ws.Range("G56").Value = 64.5
$ws.Range("J57").Value = 80.09999999999999
$ws.Range("J59").Value = 83.7
$ws.Range("D60").Value = 98
$ws.Range("D61").Value = 94.5
$ws.Range("F61").Value = 87.40000000000001
$ws.Range("G62").Value = 83.59999999999999
$ws.Range("D63").Value = 91.3
$ws.Range("E63").Value = 94.59999999999999
$ws.Range("D64").Value = 88
$ws.Range("D65").Value = 92.7
$ws.Range("D66").Value = 93.59999999999999
$ws.Range("E67").Value = 98.3
$ws.Range("D68").Value = 93.59999999999999
$ws.Range("B69").Value = 97.2
$ws.Range("G71").Value = 100.4
$ws.Range("D72").Value = 102.8
$ws.Range("G72").Value = 100.9
$ws.Range("D74").Value = 102.9
$ws.Range("H75").Value = 102.1
$ws.Range("D76").Value = 101.5
$ws.Range("D77").Value = 101.4
$ws.Range("F77").Value = 102.8
$ws.Range("B78").Value = 102.9
$ws.Range("D78").Value = 104
$ws.Range("D79").Value = 104.8
$ws.Range("G79").Value = 104.5
$ws.Range("J79").Value = 104.3
$ws.Range("D80").Value = 96.40000000000001
$ws.Range("D81").Value = 100.4
$ws.Range("H81").Value = 107.1
$ws.Range("J81").Value = 105.4
$ws.Range("D82").Value = 103.3
$ws.Range("D83").Value = 100.2
$ws.Range("G83").Value = 104.5
$ws.Range("H83").Value = 108.4
$ws.Range("B84").Value = 105.8
$ws.Range("D84").Value = 96.7
$ws.Range("G84").Value = 104.7
$ws.Range("I84").Value = 105.7
$ws.Range("F85").Value = 107.6
$ws.Range("G85").Value = 106.3
$ws.Range("H85").Value = 108.9
$ws.Range("I85").Value = 105.8
$ws.Range("J85").Value = 107
$ws.Range("C86").Value = 99.09999999999999
$ws.Range("D86").Value = 89.2
$ws.Range("I86").Value = 104.8
$ws.Range("B87").Value = 106.4
$ws.Range("D87").Value = 96.8
$ws.Range("F87").Value = 105.4
$ws.Range("G87").Value = 110.1
$ws.Range("H87").Value = 109.2
$ws.Range("J87").Value = 107.4
$ws.Range("C88").Value = 103.8
$ws.Range("D88").Value = 102.1
$ws.Range("F88").Value = 107.8
$ws.Range("I88").Value = 108
$ws.Range("C89").Value = 104.7
$ws.Range("H89").Value = 111.2
$ws.Range("I89").Value = 108.7
$ws.Range("J89").Value = 109.7
$ws.Range("D90").Value = 105.7
$ws.Range("I90").Value = 110.1
$ws.Range("B91").Value = 111.5
$ws.Range("D91").Value = 101.4
$ws.Range("E91").Value = 105.1
$ws.Range("F91").Value = 110.6
$ws.Range("G91").Value = 114.4
$ws.Range("H91").Value = 113.9
$ws.Range("I91").Value = 111.1
$ws.Range("J91").Value = 112.5
$ws.Range("B92").Value = 111
$ws.Range("C92").Value = 104.9
$ws.Range("D92").Value = 101
$ws.Range("F92").Value = 108.8
$ws.Range("G92").Value = 113.5
$ws.Range("I92").Value = 110.6
$ws.Range("B93").Value = 112.1
$ws.Range("E93").Value = 104.1
$ws.Range("F93").Value = 110.5
$ws.Range("G93").Value = 113.4
$ws.Range("H93").Value = 115.3
$ws.Range("I93").Value = 111.8
$ws.Range("J93").Value = 113
$ws.Range("B94").Value = 111.9
$ws.Range("C94").Value = 104.3
$ws.Range("D94").Value = 96.40000000000001
$ws.Range("I94").Value = 111.6
$ws.Range("B95").Value = 113.4
$ws.Range("D95").Value = 97.5
$ws.Range("F95").Value = 112.3
$ws.Range("G95").Value = 116.8
$ws.Range("H95").Value = 117.6
$ws.Range("I95").Value = 113
$ws.Range("J95").Value = 115.1
$ws.Range("B96").Value = 114
$ws.Range("C96").Value = 107.5
$ws.Range("D96").Value = 102.9
$ws.Range("F96").Value = 112.8
$ws.Range("G96").Value = 114.7
$ws.Range("I96").Value = 113.7
$ws.Range("B97").Value = 109.5
$ws.Range("F97").Value = 111
$ws.Range("G97").Value = 107.9
$ws.Range("H97").Value = 112.6
$ws.Range("I97").Value = 109.4
$ws.Range("J97").Value = 110.3
$ws.Range("C98").Value = 106.2
$ws.Range("D98").Value = 101.6
$ws.Range("E98").Value = 102.7
$ws.Range("H98").Value = 115.9
$ws.Range("J98").Value = 112.9
$ws.Range("B99").Value = 97.90000000000001
$ws.Range("C99").Value = 97.09999999999999
$ws.Range("D99").Value = 100.4
$ws.Range("F99").Value = 98.3
$ws.Range("G99").Value = 94.40000000000001
$ws.Range("H99").Value = 99.7
$ws.Range("I99").Value = 98.2
$ws.Range("J99").Value = 97.5
$ws.Range("C100").Value = 98.2
$ws.Range("D100").Value = 102.5
$ws.Range("E100").Value = 98.8
$ws.Range("F100").Value = 94.3
$ws.Range("G100").Value = 113.3
$ws.Range("J100").Value = 102.7
$ws.Range("B101").Value = 109.3
$ws.Range("E101").Value = 105.5
$ws.Range("F101").Value = 105.9
$ws.Range("G101").Value = 122.8
$ws.Range("H101").Value = 109.7
$ws.Range("I101").Value = 108.8
$ws.Range("J101").Value = 110.4
$ws.Range("C102").Value = 104.6
$ws.Range("D102").Value = 100.5
$ws.Range("H102").Value = 116.5
$ws.Range("B103").Value = 114.9
$ws.Range("C103").Value = 106.6
$ws.Range("D103").Value = 103.4
$ws.Range("F103").Value = 109.1
$ws.Range("G103").Value = 135.8
$ws.Range("H103").Value = 114.6
$ws.Range("I103").Value = 113.6
$ws.Range("J103").Value = 116.3
